$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$s.Shapes.Item(3).TextFrame.TextRange.Text = "#2 : Numerical data (SFS)"
$s.Shapes.Item(4).TextFrame.TextRange.Text = "#3 : Categorical data (SFS)"
$s.Shapes.Item(5).TextFrame.TextRange.Text = "#4 : Selected by RFE"
$s.Shapes.Item(6).TextFrame.TextRange.Text = "#5 : Numerical + Categorical (SFS)"
